# Added Week 15 simulations
# Update row 3 ("R") values on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

$offSheet = $wb.Worksheets.Item("OFF")
$offSheet.Range("B3").Value = 192
$offSheet.Range("C3").Value = 127
$offSheet.Range("D3").Value = 39
$offSheet.Range("E3").Value = 19
$offSheet.Range("F3").Value = 3
$offSheet.Range("G3").Value = 2

$defSheet = $wb.Worksheets.Item("DEF")
$defSheet.Range("B3").Value = 211
$defSheet.Range("C3").Value = 152
$defSheet.Range("D3").Value = 41
$defSheet.Range("E3").Value = 21
$defSheet.Range("F3").Value = 3
$defSheet.Range("G3").Value = 3
